$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05881766666666666
$ws.Range("H2").Value = 0.176453
$ws.Range("M2").Value = 1.116695
$ws.Range("N2").Value = 3.350085
$ws.Range("O2").Value = 0.008174214292497491
$ws.Range("P2").Value = 0.008174214292497492
$ws.Range("Q2").Value = 0.06568139427833333
$ws.Range("R2").Value = 0.591132548505
$ws.Range("S2").Value = 0.008174214292497491
$ws.Range("T2").Value = 0.008174214292497492

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05881766666666666
$ws.Range("H3").Value = 0.176453
$ws.Range("O3").Value = 0.8193429796700005
$ws.Range("P3").Value = 0.8193429796700005
$ws.Range("Q3").Value = 6.583579457451111
$ws.Range("R3").Value = 59.25221511706
$ws.Range("S3").Value = 0.8193429796700005
$ws.Range("T3").Value = 0.8193429796700005

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05881766666666666
$ws.Range("H4").Value = 0.176453
$ws.Range("O4").Value = 0.172482806037502
$ws.Range("P4").Value = 0.1724828060375021
$ws.Range("Q4").Value = 1.385932737288333
$ws.Range("R4").Value = 12.473394635595
$ws.Range("S4").Value = 0.172482806037502
$ws.Range("T4").Value = 0.1724828060375021
